# Bump the published "term" value set from version 1.0.0 to 1.1.0
# (commit message: "Added 1.1.0 of term")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Row 3: Property "Version" -> Value "1.0.0"  =>  "1.1.0"
$ws.Range("B3").Value = "1.1.0"

# Row 8: Property "Date" -> Value "2023-06-07T11:52:14+02:00"  =>  "2023-07-10T23:08:03+02:00"
$ws.Range("B8").Value = "2023-07-10T23:08:03+02:00"

# Re-assert the existing top/wrap-text alignment on the header row and the
# data rows so the alignment formatting is explicitly (re)applied, matching
# the re-saved workbook's style definitions.
$header = $ws.Range("A1:B1")
$header.VerticalAlignment = -4160
$header.WrapText = $true

$data = $ws.Range("A2:B14")
$data.VerticalAlignment = -4160
$data.WrapText = $true
